$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 476.85715
$ws.Range("I4").Value = 247.77777
$ws.Range("J4").Value = 1250
$ws.Range("K4").Value = 247.77777
$ws.Range("L4").Value = 1250
$ws.Range("M4").Value = -133.77777
$ws.Range("N4").Value = -1478
$ws.Range("H8").Value = 150.375
$ws.Range("I8").Value = 29
$ws.Range("J8").Value = 1000
$ws.Range("K8").Value = 87
$ws.Range("L8").Value = 3000
$ws.Range("M8").Value = 52
$ws.Range("N8").Value = -3278
$ws.Range("H11").Value = 3380.2144
$ws.Range("I11").Value = 3380.2144
$ws.Range("K11").Value = 3380.2144
$ws.Range("M11").Value = -3240.2144
$ws.Range("H17").Value = 3124
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3124
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9372
$ws.Range("N17").Value = -9708
$ws.Range("M17").ClearContents()
$ws.Range("H18").Value = 6460.6665
$ws.Range("I18").Value = 995.5
$ws.Range("J18").Value = 9193.25
$ws.Range("K18").Value = 995.5
$ws.Range("L18").Value = 9193.25
$ws.Range("M18").Value = -711.5
$ws.Range("N18").Value = -9761.25
$ws.Range("H26").Value = 11250
$ws.Range("J26").Value = 11250
$ws.Range("L26").Value = 11250
$ws.Range("N26").Value = -11938
$ws.Range("H30").Value = 1999
$ws.Range("J30").Value = 1999
$ws.Range("L30").Value = 5997
$ws.Range("N30").Value = -6199
$ws.Range("H37").Value = 2429.875
$ws.Range("I37").Value = 1200
$ws.Range("J37").Value = 3167.8
$ws.Range("K37").Value = 3600
$ws.Range("L37").Value = 9503.400000000001
$ws.Range("M37").Value = -3474
$ws.Range("N37").Value = -9755.400000000001
$ws.Range("H40").Value = 2750.375
$ws.Range("I40").Value = 2166.8333
$ws.Range("K40").Value = 2166.8333
$ws.Range("M40").Value = -1991.8333
$ws.Range("H43").Value = 13862.294
$ws.Range("I43").Value = 5980
$ws.Range("K43").Value = 5980
$ws.Range("M43").Value = -5911
$ws.Range("H51").Value = 23237.125
$ws.Range("I51").Value = 4919.8
$ws.Range("K51").Value = 4919.8
$ws.Range("M51").Value = -4435.8
$ws.Range("H53").Value = 469.51614
$ws.Range("I53").Value = 276.73077
$ws.Range("K53").Value = 276.73077
$ws.Range("M53").Value = 360.26923
$ws.Range("H55").Value = 1232.3846
$ws.Range("I55").Value = 971.2222
$ws.Range("J55").Value = 1820
$ws.Range("K55").Value = 971.2222
$ws.Range("L55").Value = 1820
$ws.Range("M55").Value = -757.2222
$ws.Range("N55").Value = -2248
$ws.Range("H62").Value = 100004860
$ws.Range("I62").Value = 157146780
$ws.Range("J62").Value = 6525
$ws.Range("K62").Value = 157146780
$ws.Range("L62").Value = 6525
$ws.Range("M62").Value = -157146156
$ws.Range("N62").Value = -7773
$ws.Range("H64").Value = 3720.8
$ws.Range("I64").Value = 4201.3335
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 4201.3335
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -3953.3335
$ws.Range("N64").Value = -3496
$ws.Range("H65").Value = 100004860
$ws.Range("I65").Value = 157146780
$ws.Range("J65").Value = 6525
$ws.Range("K65").Value = 785733900
$ws.Range("L65").Value = 32625
$ws.Range("M65").Value = -785730780
$ws.Range("N65").Value = -38865
$ws.Range("H67").Value = 3720.8
$ws.Range("I67").Value = 4201.3335
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 4201.3335
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -3343.3335
$ws.Range("N67").Value = -4716
$ws.Range("H76").Value = 6255854.5
$ws.Range("I76").Value = 9096291
$ws.Range("K76").Value = 9096291
$ws.Range("M76").Value = -9095976
$ws.Range("H79").Value = 6255854.5
$ws.Range("I79").Value = 9096291
$ws.Range("K79").Value = 9096291
$ws.Range("M79").Value = -9095199
$ws.Range("H80").Value = 1255
$ws.Range("J80").Value = 2315.5
$ws.Range("L80").Value = 6946.5
$ws.Range("N80").Value = -8942.5
$ws.Range("H83").Value = 1255
$ws.Range("J83").Value = 2315.5
$ws.Range("L83").Value = 20839.5
$ws.Range("N83").Value = -30823.5
$ws.Range("H96").Value = 562.1667
$ws.Range("I96").Value = 540.875
$ws.Range("J96").Value = 604.75
$ws.Range("K96").Value = 1622.625
$ws.Range("L96").Value = 1814.25
$ws.Range("M96").Value = -249.625
$ws.Range("N96").Value = -4560.25
$ws.Range("H97").Value = 910.5
$ws.Range("J97").Value = 910.5
$ws.Range("L97").Value = 2731.5
$ws.Range("N97").Value = -3723.5
$ws.Range("H98").Value = 7625.3335
$ws.Range("I98").Value = 8243.546
$ws.Range("K98").Value = 8243.546
$ws.Range("M98").Value = -6745.546
$ws.Range("H112").Value = 2739
$ws.Range("J112").Value = 2739
$ws.Range("L112").Value = 8217
$ws.Range("N112").Value = -10433
$ws.Range("H117").Value = 80666.336
$ws.Range("J117").Value = 80666.336
$ws.Range("L117").Value = 80666.336
$ws.Range("N117").Value = -89844.336
$ws.Range("H122").Value = 7625.3335
$ws.Range("I122").Value = 8243.546
$ws.Range("K122").Value = 24730.638
$ws.Range("M122").Value = -22280.638
$ws.Range("H125").Value = 2867.077
$ws.Range("I125").Value = 3027.625
$ws.Range("K125").Value = 27248.625
$ws.Range("M125").Value = -24788.625
$ws.Range("H137").Value = 7504.385
$ws.Range("J137").Value = 7351.4443
$ws.Range("L137").Value = 22054.3329
$ws.Range("N137").Value = -27154.3329
$ws.Range("H138").Value = 5570.84
$ws.Range("J138").Value = 6183.427
$ws.Range("L138").Value = 18550.281
$ws.Range("N138").Value = -28830.281
$ws.Range("H141").Value = 1659.5151
$ws.Range("I141").Value = 1525.4
$ws.Range("K141").Value = 4576.200000000001
$ws.Range("M141").Value = 603.7999999999993

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 1972
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H32").Value = 32391.902
$ws.Range("I32").Value = 31319.064
$ws.Range("J32").Value = 44997.75
$ws.Range("K32").Value = 31319.064
$ws.Range("L32").Value = 44997.75
$ws.Range("M32").Value = -31032.064
$ws.Range("N32").Value = -45571.75
$ws.Range("H45").Value = 1690.7354
$ws.Range("I45").Value = 1578.174
$ws.Range("J45").Value = 1926.091
$ws.Range("K45").Value = 1578.174
$ws.Range("L45").Value = 1926.091
$ws.Range("M45").Value = -1201.174
$ws.Range("N45").Value = -2680.091
$ws.Range("H61").Value = 9113.736999999999
$ws.Range("I61").Value = 2059.4
$ws.Range("K61").Value = 2059.4
$ws.Range("M61").Value = -1847.4
$ws.Range("H74").Value = 457365
$ws.Range("I74").Value = 716742.6
$ws.Range("K74").Value = 716742.6
$ws.Range("M74").Value = -715868.6
$ws.Range("H77").Value = 457365
$ws.Range("I77").Value = 716742.6
$ws.Range("K77").Value = 3583713
$ws.Range("M77").Value = -3579345
$ws.Range("H102").Value = 7500
$ws.Range("I102").Value = 7500
$ws.Range("K102").Value = 7500
$ws.Range("M102").Value = -5878
$ws.Range("H103").Value = 46362
$ws.Range("J103").Value = 46362
$ws.Range("L103").Value = 46362
$ws.Range("N103").Value = -48706
$ws.Range("H110").Value = 6946785
$ws.Range("I110").Value = 10871229
$ws.Range("K110").Value = 10871229
$ws.Range("M110").Value = -10869184
$ws.Range("H122").Value = 3521.5557
$ws.Range("I122").Value = 2166.1
$ws.Range("J122").Value = 5215.875
$ws.Range("K122").Value = 6498.299999999999
$ws.Range("L122").Value = 15647.625
$ws.Range("M122").Value = -4048.299999999999
$ws.Range("N122").Value = -20547.625
$ws.Range("H132").Value = 9082.672
$ws.Range("I132").Value = 3853.4468
$ws.Range("J132").Value = 23539.941
$ws.Range("K132").Value = 11560.3404
$ws.Range("L132").Value = 70619.823
$ws.Range("M132").Value = -9030.340400000001
$ws.Range("N132").Value = -75679.823
$ws.Range("H136").Value = 9113.736999999999
$ws.Range("I136").Value = 2059.4
$ws.Range("K136").Value = 6178.200000000001
$ws.Range("M136").Value = -3628.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 6000
$ws.Range("J15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("N15").Value = -10454
$ws.Range("H24").Value = 1015.6667
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()
$ws.Range("H103").Value = 49000
$ws.Range("J103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H104").Value = 85000
$ws.Range("J104").Value = 85000
$ws.Range("L104").Value = 85000
$ws.Range("N104").Value = -91988
$ws.Range("H105").Value = 28581142
$ws.Range("I105").Value = 41679040
$ws.Range("K105").Value = 41679040
$ws.Range("M105").Value = -41677293
$ws.Range("H106").Value = 71100
$ws.Range("J106").Value = 71100
$ws.Range("L106").Value = 71100
$ws.Range("N106").Value = -73624
$ws.Range("H107").Value = 2444.6667
$ws.Range("I107").Value = 1782.64
$ws.Range("K107").Value = 1782.64
$ws.Range("M107").Value = 137.3599999999999
$ws.Range("H134").Value = 5733.12
$ws.Range("I134").Value = 4817.825
$ws.Range("K134").Value = 14453.475
$ws.Range("M134").Value = -11918.475

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 914500
$ws.Range("I3").Value = 2254750
$ws.Range("J3").Value = 21000
$ws.Range("K3").Value = 2254750
$ws.Range("L3").Value = 21000
$ws.Range("M3").Value = -2254637
$ws.Range("N3").Value = -21226
$ws.Range("H14").Value = 3661.7
$ws.Range("J14").Value = 2423.4
$ws.Range("L14").Value = 2423.4
$ws.Range("N14").Value = -2763.4
$ws.Range("H17").Value = 30000
$ws.Range("J17").Value = 30000
$ws.Range("L17").Value = 30000
$ws.Range("N17").Value = -30348
$ws.Range("H22").Value = 519.9167
$ws.Range("I22").Value = 430.86365
$ws.Range("K22").Value = 430.86365
$ws.Range("M22").Value = -80.86365000000001
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
$ws.Range("H31").Value = 8122.852
$ws.Range("J31").Value = 9862.210999999999
$ws.Range("L31").Value = 9862.210999999999
$ws.Range("N31").Value = -10452.211
$ws.Range("H34").Value = 8122.852
$ws.Range("J34").Value = 9862.210999999999
$ws.Range("L34").Value = 9862.210999999999
$ws.Range("N34").Value = -10266.211
$ws.Range("H59").Value = 49728.145
$ws.Range("I59").Value = 33619.6
$ws.Range("J59").Value = 89999.5
$ws.Range("K59").Value = 33619.6
$ws.Range("L59").Value = 89999.5
$ws.Range("M59").Value = -32474.6
$ws.Range("N59").Value = -92289.5
$ws.Range("H103").Value = 4000
$ws.Range("I103").Value = 4000
$ws.Range("K103").Value = 4000
$ws.Range("M103").Value = -2828
$ws.Range("H132").Value = 24394.107
$ws.Range("I132").Value = 4464.684
$ws.Range("K132").Value = 13394.052
$ws.Range("M132").Value = -10864.052
$ws.Range("H134").Value = 3697.2122
$ws.Range("I134").Value = 2631.3076
$ws.Range("J134").Value = 7656.2856
$ws.Range("K134").Value = 7893.9228
$ws.Range("L134").Value = 22968.8568
$ws.Range("M134").Value = -5358.9228
$ws.Range("N134").Value = -28038.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4148.294
$ws.Range("I2").Value = 54.88889
$ws.Range("J2").Value = 8753.375
$ws.Range("K2").Value = 329.33334
$ws.Range("L2").Value = 52520.25
$ws.Range("M2").Value = -216.33334
$ws.Range("N2").Value = -52746.25
$ws.Range("H4").Value = 4263067.5
$ws.Range("I4").Value = 2779737.2
$ws.Range("J4").Value = 8713058
$ws.Range("K4").Value = 8339211.600000001
$ws.Range("L4").Value = 26139174
$ws.Range("M4").Value = -8339099.600000001
$ws.Range("N4").Value = -26139398
$ws.Range("H7").Value = 100
$ws.Range("J7").Value = 100
$ws.Range("L7").Value = 300
$ws.Range("N7").Value = -524
$ws.Range("H10").Value = 558.875
$ws.Range("I10").Value = 294
$ws.Range("J10").Value = 1000.3333
$ws.Range("K10").Value = 882
$ws.Range("L10").Value = 3000.9999
$ws.Range("M10").Value = -743
$ws.Range("N10").Value = -3278.9999
$ws.Range("H12").Value = 152.52942
$ws.Range("J12").Value = 183.21428
$ws.Range("L12").Value = 549.64284
$ws.Range("N12").Value = -895.64284
$ws.Range("H38").Value = 83333370
$ws.Range("I38").Value = 250000050
$ws.Range("K38").Value = 750000150
$ws.Range("M38").Value = -749999803
$ws.Range("H45").Value = 3584.75
$ws.Range("J45").Value = 4113
$ws.Range("L45").Value = 12339
$ws.Range("N45").Value = -13403
$ws.Range("H68").Value = 1727.4546
$ws.Range("J68").Value = 2695.3333
$ws.Range("L68").Value = 8085.999899999999
$ws.Range("N68").Value = -9707.999899999999
$ws.Range("H71").Value = 1727.4546
$ws.Range("J71").Value = 2695.3333
$ws.Range("L71").Value = 24257.9997
$ws.Range("N71").Value = -32369.9997
$ws.Range("H94").Value = 4199.8
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352
$ws.Range("H113").Value = 3690.625
$ws.Range("J113").Value = 4161.5
$ws.Range("L113").Value = 12484.5
$ws.Range("N113").Value = -16824.5
$ws.Range("H131").Value = 25662032
$ws.Range("I131").Value = 83334584
$ws.Range("J131").Value = 29784.777
$ws.Range("K131").Value = 250003752
$ws.Range("L131").Value = 89354.33099999999
$ws.Range("M131").Value = -249998712
$ws.Range("N131").Value = -99434.33099999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H20").Value = 106000
$ws.Range("J20").Value = 12000
$ws.Range("L20").Value = 12000
$ws.Range("N20").Value = -12490
$ws.Range("H24").Value = 134335.67
$ws.Range("J24").Value = 3007
$ws.Range("L24").Value = 3007
$ws.Range("N24").Value = -3353
$ws.Range("H26").Value = 46999
$ws.Range("J26").Value = 46999
$ws.Range("L26").Value = 46999
$ws.Range("N26").Value = -47559
$ws.Range("H50").Value = 46999
$ws.Range("J50").Value = 46999
$ws.Range("L50").Value = 46999
$ws.Range("N50").Value = -47995
$ws.Range("H80").Value = 7001.6665
$ws.Range("I80").Value = 7001.6665
$ws.Range("K80").Value = 7001.6665
$ws.Range("M80").Value = -6003.6665
$ws.Range("H83").Value = 7001.6665
$ws.Range("I83").Value = 7001.6665
$ws.Range("K83").Value = 35008.3325
$ws.Range("M83").Value = -30016.3325
$ws.Range("H113").Value = 980.5
$ws.Range("I113").Value = 980.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 980.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1189.5
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2299.5
$ws.Range("I122").Value = 2024.375
$ws.Range("K122").Value = 6073.125
$ws.Range("M122").Value = -3623.125
$ws.Range("H126").Value = 3495.4443
$ws.Range("I126").Value = 2779.8572
$ws.Range("K126").Value = 8339.571599999999
$ws.Range("M126").Value = -5869.571599999999
$ws.Range("H132").Value = 6647.5
$ws.Range("I132").Value = 4860.273
$ws.Range("K132").Value = 14580.819
$ws.Range("M132").Value = -12050.819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4187.875
$ws.Range("I7").Value = 4093.2856
$ws.Range("K7").Value = 4093.2856
$ws.Range("M7").Value = -3981.2856
$ws.Range("H22").Value = 3558.2727
$ws.Range("I22").Value = 2139.3076
$ws.Range("J22").Value = 4480.6
$ws.Range("K22").Value = 2139.3076
$ws.Range("L22").Value = 4480.6
$ws.Range("M22").Value = -1844.3076
$ws.Range("N22").Value = -5070.6
$ws.Range("H27").Value = 3558.2727
$ws.Range("I27").Value = 2139.3076
$ws.Range("J27").Value = 4480.6
$ws.Range("K27").Value = 2139.3076
$ws.Range("L27").Value = 4480.6
$ws.Range("M27").Value = -2032.3076
$ws.Range("N27").Value = -4694.6
$ws.Range("H46").Value = 5034.913
$ws.Range("I46").Value = 1521.6666
$ws.Range("K46").Value = 1521.6666
$ws.Range("M46").Value = -1333.6666
$ws.Range("H55").Value = 644.0526
$ws.Range("I55").Value = 670.36365
$ws.Range("K55").Value = 670.36365
$ws.Range("M55").Value = -497.36365
$ws.Range("H61").Value = 3613.353
$ws.Range("I61").Value = 2602.7273
$ws.Range("K61").Value = 2602.7273
$ws.Range("M61").Value = -2400.7273
$ws.Range("H113").Value = 3613.353
$ws.Range("I113").Value = 2602.7273
$ws.Range("K113").Value = 2602.7273
$ws.Range("M113").Value = -432.7273
$ws.Range("H122").Value = 33337854
$ws.Range("I122").Value = 50004412
$ws.Range("J122").Value = 4734.8
$ws.Range("K122").Value = 150013236
$ws.Range("L122").Value = 14204.4
$ws.Range("M122").Value = -150010786
$ws.Range("N122").Value = -19104.4
$ws.Range("H126").Value = 4187.875
$ws.Range("I126").Value = 4093.2856
$ws.Range("K126").Value = 12279.8568
$ws.Range("M126").Value = -9809.856800000001
$ws.Range("H136").Value = 7888.92
$ws.Range("I136").Value = 7081.857
$ws.Range("K136").Value = 21245.571
$ws.Range("M136").Value = -18695.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 287200.3
$ws.Range("J2").Value = 341249.88
$ws.Range("L2").Value = 341249.88
$ws.Range("N2").Value = -341473.88
$ws.Range("H31").Value = 22999.5
$ws.Range("J31").Value = 22999.5
$ws.Range("L31").Value = 22999.5
$ws.Range("N31").Value = -23695.5
$ws.Range("H96").Value = 5958.3335
$ws.Range("I96").Value = 5017.3335
$ws.Range("J96").Value = 6899.3335
$ws.Range("K96").Value = 5017.3335
$ws.Range("L96").Value = 6899.3335
$ws.Range("M96").Value = -3644.3335
$ws.Range("N96").Value = -9645.333500000001
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H107").Value = 2123.1667
$ws.Range("J107").Value = 2123.1667
$ws.Range("L107").Value = 6369.500100000001
$ws.Range("N107").Value = -10209.5001
$ws.Range("H126").Value = 1178.875
$ws.Range("I126").Value = 990.2857
$ws.Range("J126").Value = 2499
$ws.Range("K126").Value = 2970.8571
$ws.Range("L126").Value = 7497
$ws.Range("M126").Value = -500.8571000000002
$ws.Range("N126").Value = -12437
$ws.Range("H132").Value = 3787.0625
$ws.Range("I132").Value = 2985.7273
$ws.Range("J132").Value = 5550
$ws.Range("K132").Value = 8957.1819
$ws.Range("L132").Value = 16650
$ws.Range("M132").Value = -6427.1819
$ws.Range("N132").Value = -21710

